$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.00%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07518"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.593"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.16%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9270"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.04%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1182"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.98%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1835"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.34%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08905"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.50%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04125"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.63%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.66%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001283"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005785"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.80%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.337"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.18%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.06%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3332"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.34%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.324"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.26%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.83%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.97%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04103"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.55%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.44%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003890"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.06%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.34%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02393"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.76%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05225"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006684"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "16.79%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007797"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1325"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.07%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007401"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007115"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.88%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2992"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.82%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006579"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.19%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "13.07%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
